$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.01441373773993604
$ws.Range("D2").Value = 0.04580567640922517
$ws.Range("E2").Value = 0.06496870325029391
$ws.Range("F2").Value = 1.47758319985725
$ws.Range("G2").Value = 0.00246675029531099
$ws.Range("I2").Value = 1.133028412353909
$ws.Range("K2").Value = 1.592004394650132
$ws.Range("M2").Value = 0.4553131149134231
$ws.Range("N2").Value = 1.469867471075457

$ws.Range("C3").Value = 0.0142020115517667
$ws.Range("D3").Value = 0.04650678353788162
$ws.Range("E3").Value = 0.06050874881020718
$ws.Range("F3").Value = 1.445795482841334
$ws.Range("G3").Value = 0.002471879432849382
$ws.Range("I3").Value = 1.109565171040089
$ws.Range("K3").Value = 1.438226604960391
$ws.Range("M3").Value = 0.4144391994493049
$ws.Range("N3").Value = 1.489657874149604

$ws.Range("C4").Value = 0.01407048987438841
$ws.Range("D4").Value = 0.04695873728501532
$ws.Range("E4").Value = 0.05782092484009027
$ws.Range("F4").Value = 1.427364385868913
$ws.Range("G4").Value = 0.00247519242152506
$ws.Range("I4").Value = 1.095994222341062
$ws.Range("K4").Value = 1.344468734700058
$ws.Range("M4").Value = 0.3895773071537789
$ws.Range("N4").Value = 1.502425090868535

$ws.Range("C5").Value = 0.01401651720777863
$ws.Range("D5").Value = 0.04714828660392367
$ws.Range("E5").Value = 0.0567381089546366
$ws.Range("F5").Value = 1.420125004416448
$ws.Range("G5").Value = 0.002476583794069935
$ws.Range("I5").Value = 1.090672653471117
$ws.Range("K5").Value = 1.306426183474059
$ws.Range("M5").Value = 0.3795040764754489
$ws.Range("N5").Value = 1.507782254415943

$ws.Range("C6").Value = 0.01400753252353582
$ws.Range("D6").Value = 0.04718008502018911
$ws.Range("E6").Value = 0.05655905610478129
$ws.Range("F6").Value = 1.418939252182057
$ws.Range("G6").Value = 0.002476817329321394
$ws.Range("I6").Value = 1.089801575203737
$ws.Range("K6").Value = 1.300119127206358
$ws.Range("M6").Value = 0.3778349138669839
$ws.Range("N6").Value = 1.50868112121838

$ws.Range("C7").Value = 0.01406976349723976
$ws.Range("D7").Value = 0.04696127188530497
$ws.Range("E7").Value = 0.0578062713205334
$ws.Range("F7").Value = 1.427265656398518
$ws.Range("G7").Value = 0.002475211018883705
$ws.Range("I7").Value = 1.095921610642066
$ws.Range("K7").Value = 1.343955015938263
$ws.Range("M7").Value = 0.3894412215542857
$ws.Range("N7").Value = 1.502496714807503

$ws.Range("C8").Value = 0.01434105367398786
$ws.Range("D8").Value = 0.0460429469779573
$ws.Range("E8").Value = 0.06342026661659972
$ws.Range("F8").Value = 1.466396140827641
$ws.Range("G8").Value = 0.002468484938114076
$ws.Range("I8").Value = 1.124764000755448
$ws.Range("K8").Value = 1.538843084454811
$ws.Range("M8").Value = 0.4411705078414556
$ws.Range("N8").Value = 1.4765631251956

$ws.Range("C9").Value = 0.01486075376038798
$ws.Range("D9").Value = 0.04441356500788984
$ws.Range("E9").Value = 0.07484189357784743
$ws.Range("F9").Value = 1.551841610189911
$ws.Range("G9").Value = 0.002456587205193674
$ws.Range("I9").Value = 1.188022109002816
$ws.Range("K9").Value = 1.926385209331841
$ws.Range("M9").Value = 0.54451971749981
$ws.Range("N9").Value = 1.430612586094476

$ws.Range("C10").Value = 0.01523481248904801
$ws.Range("D10").Value = 0.04332245137903357
$ws.Range("E10").Value = 0.08350113513529323
$ws.Range("F10").Value = 1.620060158376887
$ws.Range("G10").Value = 0.002448624349091655
$ws.Range("I10").Value = 1.238682809170783
$ws.Range("K10").Value = 2.214570945234584
$ws.Range("M10").Value = 0.6216839780642545
$ws.Range("N10").Value = 1.399867879237409

$ws.Range("C11").Value = 0.01540324411847749
$ws.Range("D11").Value = 0.0428494172966527
$ws.Range("E11").Value = 0.08750218806536481
$ws.Range("F11").Value = 1.652305406701231
$ws.Range("G11").Value = 0.00244516887813861
$ws.Range("I11").Value = 1.262660963234865
$ws.Range("K11").Value = 2.346467888578218
$ws.Range("M11").Value = 0.6570717257745429
$ws.Range("N11").Value = 1.386541471194644

$ws.Range("C12").Value = 0.01546677129669405
$ws.Range("D12").Value = 0.04267367427107871
$ws.Range("E12").Value = 0.08902646526360058
$ws.Range("F12").Value = 1.664692406895767
$ws.Range("G12").Value = 0.002443884225355515
$ws.Range("I12").Value = 1.271876668749229
$ws.Range("K12").Value = 2.39653173981759
$ws.Range("M12").Value = 0.6705143168051535
$ws.Range("N12").Value = 1.381590450238619

$ws.Range("C13").Value = 0.01545310097477781
$ws.Range("D13").Value = 0.04271137263016733
$ws.Range("E13").Value = 0.08869777307162252
$ws.Range("F13").Value = 1.662016764520502
$ws.Range("G13").Value = 0.002444159839199767
$ws.Range("I13").Value = 1.269885840260912
$ws.Range("K13").Value = 2.385744337203448
$ws.Range("M13").Value = 0.6676173295850703
$ws.Range("N13").Value = 1.382652489138152

$ws.Range("C14").Value = 0.01540847566410619
$ws.Range("D14").Value = 0.04283489079381653
$ws.Range("E14").Value = 0.08762740607318165
$ws.Range("F14").Value = 1.653320944731888
$ws.Range("G14").Value = 0.002445062711618846
$ws.Range("I14").Value = 1.263416415859638
$ws.Range("K14").Value = 2.350584306506278
$ws.Range("M14").Value = 0.6581768079841055
$ws.Range("N14").Value = 1.386132232490404

$ws.Range("C15").Value = 0.01538110810364302
$ws.Range("D15").Value = 0.04291099097165763
$ws.Range("E15").Value = 0.08697297678360627
$ws.Range("F15").Value = 1.64801754246659
$ws.Range("G15").Value = 0.002445618849910964
$ws.Range("I15").Value = 1.259471427326147
$ws.Range("K15").Value = 2.329063123906678
$ws.Range("M15").Value = 0.6523997176469152
$ws.Range("N15").Value = 1.388276116624414

$ws.Range("C16").Value = 0.01522376980637574
$ws.Range("D16").Value = 0.04335383718042429
$ws.Range("E16").Value = 0.08324092516399872
$ws.Range("F16").Value = 1.61797742868896
$ws.Range("G16").Value = 0.002448853518524293
$ws.Range("I16").Value = 1.237134680344724
$ws.Range("K16").Value = 2.205967461852254
$ws.Range("M16").Value = 0.6193771312889282
$ws.Range("N16").Value = 1.400752110363197

$ws.Range("C17").Value = 0.01512680089311402
$ws.Range("D17").Value = 0.04363150370184776
$ws.Range("E17").Value = 0.08096747645513602
$ws.Range("F17").Value = 1.599860772629995
$ws.Range("G17").Value = 0.002450880525130761
$ws.Range("I17").Value = 1.223671758024196
$ws.Range("K17").Value = 2.13065875012677
$ws.Range("M17").Value = 0.5991926233202207
$ws.Range("N17").Value = 1.408575043232855

$ws.Range("C18").Value = 0.01507086446238404
$ws.Range("D18").Value = 0.04379340310786617
$ws.Range("E18").Value = 0.07966566131281638
$ws.Range("F18").Value = 1.589554530416933
$ws.Range("G18").Value = 0.002452062120473346
$ws.Range("I18").Value = 1.216015899483182
$ws.Range("K18").Value = 2.087418309798522
$ws.Range("M18").Value = 0.5876098254698121
$ws.Range("N18").Value = 1.413136629545793

$ws.Range("C19").Value = 0.01505189761410364
$ws.Range("D19").Value = 0.04384859516148865
$ws.Range("E19").Value = 0.07922587946215742
$ws.Range("F19").Value = 1.586084526413032
$ws.Range("G19").Value = 0.002452464891539311
$ws.Range("I19").Value = 1.213438763443548
$ws.Range("K19").Value = 2.072790674764178
$ws.Range("M19").Value = 0.5836926617947569
$ws.Range("N19").Value = 1.414691743090852

$ws.Range("C20").Value = 0.01513714025568902
$ws.Range("D20").Value = 0.04360171843726057
$ws.Range("E20").Value = 0.08120888550651273
$ws.Range("F20").Value = 1.601777510070704
$ws.Range("G20").Value = 0.002450663121520203
$ws.Range("I20").Value = 1.225095825061771
$ws.Range("K20").Value = 2.138667697769449
$ws.Range("M20").Value = 0.6013385181374531
$ws.Range("N20").Value = 1.407735853492252

$ws.Range("C21").Value = 0.01542159014307742
$ws.Range("D21").Value = 0.04279851840014715
$ws.Range("E21").Value = 0.08794154796454023
$ws.Range("F21").Value = 1.655870312816774
$ws.Range("G21").Value = 0.002444796869614441
$ws.Range("I21").Value = 1.265312946709116
$ws.Range("K21").Value = 2.360908456353229
$ws.Range("M21").Value = 0.6609485696767052
$ws.Range("N21").Value = 1.385107553833315

$ws.Range("C22").Value = 0.0156060103319291
$ws.Range("D22").Value = 0.04229332458296753
$ws.Range("E22").Value = 0.09239527812926696
$ws.Range("F22").Value = 1.692252417405797
$ws.Range("G22").Value = 0.002441101940470869
$ws.Range("I22").Value = 1.29238888935356
$ws.Range("K22").Value = 2.50684071635942
$ws.Range("M22").Value = 0.7001525743595636
$ws.Range("N22").Value = 1.370874987898446

$ws.Range("C23").Value = 0.01550771939215778
$ws.Range("D23").Value = 0.0425611388792877
$ws.Range("E23").Value = 0.09001325183675135
$ws.Range("F23").Value = 1.672739718776114
$ws.Range("G23").Value = 0.002443061319925856
$ws.Range("I23").Value = 1.277864962176452
$ws.Range("K23").Value = 2.428890453517795
$ws.Range("M23").Value = 0.6792058750116325
$ws.Range("N23").Value = 1.378420082645164

$ws.Range("C24").Value = 0.01513246641317778
$ws.Range("D24").Value = 0.04361517730154496
$ws.Range("E24").Value = 0.08109972823244505
$ws.Range("F24").Value = 1.600910612896968
$ws.Range("G24").Value = 0.002450761358918249
$ws.Range("I24").Value = 1.224451742353423
$ws.Range("K24").Value = 2.135046679716936
$ws.Range("M24").Value = 0.6003682920581497
$ws.Range("N24").Value = 1.40811505162764

$ws.Range("C25").Value = 0.01472150609872713
$ws.Range("D25").Value = 0.04483585136564372
$ws.Range("E25").Value = 0.07170620205537404
$ws.Range("F25").Value = 1.527780308802789
$ws.Range("G25").Value = 0.00245966849174437
$ws.Range("I25").Value = 1.17018159740968
$ws.Range("K25").Value = 1.820951110020474
$ws.Range("M25").Value = 0.5163495234981781
$ws.Range("N25").Value = 1.442515762178108
